$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date serial from 45233 to 45243 for rows 2-8
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45243
}
